$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Rename sheets: Hoja1 -> Partes, Hoja2 -> Pedido
$ws1.Name = "Partes"
$ws2.Name = "Pedido"

# Remove the "Notas" column (column B) from the Pedido sheet, including
# its header and the sample note value. This also prunes the now-unused
# shared strings ("Notas" / "Esta es una nota de prueba").
$ws2.Columns.Item(2).Delete()

# Update selections on each sheet. Selecting on Pedido first, then on
# Partes, leaves Partes as the final active sheet/tab (matching the
# target workbook view).
$ws2.Range("E21").Select()
$ws1.Range("A29").Select()
